$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.533.41"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "1.812.98"
$ws.Range("E3").Value = "  +0.02%  "

$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").Value = "'305.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.87%  "

$ws.Range("D7").Value = "'0.4545"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.19%  "

$ws.Range("D8").Value = "'0.3590"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.04%  "

$ws.Range("D9").Value = "'46.28"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.61%  "

$ws.Range("D10").Value = "'0.07107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.44%  "

$ws.Range("D11").Value = "'0.8916"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.30%  "

$ws.Range("D12").Value = "'0.07726"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.51%  "

$ws.Range("D13").Value = "'19.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.37%  "

$ws.Range("D14").Value = "1.777.59"
$ws.Range("E14").Value = "  -2.01%  "

$ws.Range("D15").Value = "'5.256"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "

$ws.Range("D16").Value = "'6.299"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").Value = "'85.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.96%  "

$ws.Range("E18").Value = "  -0.29%  "

$ws.Range("D19").Value = "'0.000008543"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "

$ws.Range("E20").Value = "  -0.25%  "

$ws.Range("D21").Value = "26.575.96"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "'14.13"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.84%  "

$ws.Range("D23").Value = "'4.953"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.08%  "

$ws.Range("E24").Value = "  +0.25%  "

$ws.Range("D25").Value = "'1.925"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.17%  "

$ws.Range("D26").Value = "'152.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "

$ws.Range("D27").Value = "'17.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "'2.018"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.53%  "

$ws.Range("D29").Value = "'112.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "

$ws.Range("D30").Value = "'4.817"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.84%  "

$ws.Range("D31").Value = "'0.08709"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.30%  "

$ws.Range("D32").Value = "'3.134"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.40%  "

$ws.Range("D33").Value = "'0.7418"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.11%  "

$ws.Range("D34").Value = "'4.426"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.19%  "

$ws.Range("D35").Value = "'2.705"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.71%  "

$ws.Range("D36").Value = "'1.108"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.02%  "

$ws.Range("D37").Value = "'1.071"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.11%  "

$ws.Range("D38").Value = "'0.01935"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.93%  "

$ws.Range("D39").Value = "'2.913"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("E40").Value = "  -0.73%  "

$ws.Range("D41").Value = "'0.5085"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("D42").Value = "'6.799"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("E43").Value = "  -3.45%  "

$ws.Range("D44").Value = "'8.017"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.99%  "

$ws.Range("D45").Value = "'0.4682"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("E46").Value = "  -0.31%  "

$ws.Range("D47").Value = "'9.971"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "

$ws.Range("D48").Value = "'98.97"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.04%  "

$ws.Range("D49").Value = "'1.562"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.01%  "

$ws.Range("D50").Value = "'0.05993"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.16%  "

$ws.Range("D51").Value = "'63.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.11%  "
